$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": update dates and forecast numbers ---
$ws1.Range("B2").Value = "'2025-02-02"
$ws1.Range("B3").Value = "'2025-02-09"
$ws1.Range("H3").Value = 29
$ws1.Range("B4").Value = "'2025-02-16"
$ws1.Range("H4").Value = 29
$ws1.Range("B5").Value = "'2025-02-23"
$ws1.Range("F5").Value = 18
$ws1.Range("G5").Value = 24
$ws1.Range("B6").Value = "'2025-03-02"
$ws1.Range("F6").Value = 16
$ws1.Range("G6").Value = 23
$ws1.Range("H6").Value = 35
$ws1.Range("B7").Value = "'2025-03-09"
$ws1.Range("F7").Value = 16
$ws1.Range("G7").Value = 23
$ws1.Range("B8").Value = "'2025-03-16"
$ws1.Range("D8").Value = 3
$ws1.Range("E8").Value = 15
$ws1.Range("F8").Value = 15
$ws1.Range("G8").Value = 24
$ws1.Range("H8").Value = 42
$ws1.Range("B9").Value = "'2025-03-23"
$ws1.Range("D9").Value = 3
$ws1.Range("E9").Value = 14
$ws1.Range("F9").Value = 14
$ws1.Range("G9").Value = 23
$ws1.Range("H9").Value = 42
$ws1.Range("B10").Value = "'2025-03-30"
$ws1.Range("E10").Value = 13
$ws1.Range("F10").Value = 14
$ws1.Range("G10").Value = 22
$ws1.Range("H10").Value = 37
$ws1.Range("B11").Value = "'2025-04-06"
$ws1.Range("E11").Value = 14
$ws1.Range("F11").Value = 14
$ws1.Range("G11").Value = 24
$ws1.Range("H11").Value = 43
$ws1.Range("B12").Value = "'2025-04-13"
$ws1.Range("E12").Value = 14
$ws1.Range("F12").Value = 15
$ws1.Range("G12").Value = 24
$ws1.Range("H12").Value = 41
$ws1.Range("B13").Value = "'2025-04-20"
$ws1.Range("E13").Value = 14
$ws1.Range("F13").Value = 15
$ws1.Range("G13").Value = 24
$ws1.Range("H13").Value = 41
$ws1.Range("B14").Value = "'2025-04-27"
$ws1.Range("E14").Value = 14
$ws1.Range("F14").Value = 15
$ws1.Range("G14").Value = 24
$ws1.Range("H14").Value = 40
$ws1.Range("B15").Value = "'2025-05-04"
$ws1.Range("E15").Value = 14
$ws1.Range("F15").Value = 14
$ws1.Range("G15").Value = 23
$ws1.Range("H15").Value = 39
$ws1.Range("B16").Value = "'2025-05-11"
$ws1.Range("D16").Value = 3
$ws1.Range("E16").Value = 15
$ws1.Range("F16").Value = 15
$ws1.Range("G16").Value = 24
$ws1.Range("H16").Value = 43
$ws1.Range("B17").Value = "'2025-05-18"
$ws1.Range("E17").Value = 14
$ws1.Range("F17").Value = 13
$ws1.Range("G17").Value = 22
$ws1.Range("H17").Value = 40

# --- Sheet "Summary": update metric values ---
$ws2.Range("B2").Value = "2023-01-01 to 2025-01-26"
$ws2.Range("B4").Value = "'91"
$ws2.Range("B6").Value = "'30"
$ws2.Range("B8").Value = "2913 units"
$ws2.Range("B9").Value = "'53"
$ws2.Range("B10").Value = "'26"
$ws2.Range("B11").Value = "'12"
$ws2.Range("B13").Value = "'2025-02-23"
$ws2.Range("B15").Value = "'2025-02-02"
